# Removal of unused methods (machines) from the scenario grid.
# Machine 9, Machine 10 and Machine 11 (rows 12-14) are no longer used,
# so delete those three entire rows; everything below shifts up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(12).Resize(3).EntireRow.Delete()

# Update the active selection to match the saved state after the edit.
$ws.Range("E34").Select()
